$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.204.47'

$ws.Range('D3').Value = '1.602.15'
$ws.Range('E3').Value = '  +0.28%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.07%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '304.96'
$ws.Range('E6').Value = '  +0.94%  '

$ws.Range('E7').Value = '  -0.40%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '52.87'
$ws.Range('E8').Value = '  +3.86%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3595'
$ws.Range('E9').Value = '  -0.99%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.255'
$ws.Range('E10').Value = '  +0.53%  '

$ws.Range('E11').Value = '  +0.05%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08127'
$ws.Range('E12').Value = '  -0.06%  '

$ws.Range('E13').Value = '  +2.03%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.587'
$ws.Range('E14').Value = '  +0.38%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.335'
$ws.Range('E15').Value = '  -0.27%  '

$ws.Range('E16').Value = '  +0.02%  '

$ws.Range('D17').Value = '1.601.67'
$ws.Range('E17').Value = '  +0.19%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '93.87'
$ws.Range('E18').Value = '  +1.73%  '

$ws.Range('E19').Value = '  +1.09%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.13'
$ws.Range('E20').Value = '  -0.04%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.517'
$ws.Range('E21').Value = '  +0.33%  '

$ws.Range('E22').Value = '  +0.21%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.87'
$ws.Range('E23').Value = '  -0.86%  '

$ws.Range('D24').Value = '23.195.76'
$ws.Range('E24').Value = '  +0.73%  '

$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.047'
$ws.Range('E25').Value = '  +9.24%  '

$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.419'
$ws.Range('E26').Value = '  +1.91%  '

$ws.Range('E27').Value = '  +0.24%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '150.52'
$ws.Range('E28').Value = '  +0.91%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.258'
$ws.Range('E29').Value = '  -0.07%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '134.73'
$ws.Range('E30').Value = '  -0.24%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.408'
$ws.Range('E31').Value = '  +2.10%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.725'
$ws.Range('E32').Value = '  -0.16%  '

$ws.Range('D33').Value = '1.782.60'
$ws.Range('E33').Value = '  +0.38%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9478'
$ws.Range('E34').Value = '  -0.90%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02759'
$ws.Range('E35').Value = '  +1.95%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.07391'
$ws.Range('E36').Value = '  -1.66%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '10.26'
$ws.Range('E37').Value = '  +1.32%  '

$ws.Range('E38').Value = '  -0.18%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.087'
$ws.Range('E39').Value = '  -1.60%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.08736'
$ws.Range('E40').Value = '  -0.87%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.396'
$ws.Range('E41').Value = '  +2.86%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7077'
$ws.Range('E42').Value = '  +0.85%  '

$ws.Range('E43').Value = '  +0.24%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.82'
$ws.Range('E44').Value = '  +3.48%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6504'
$ws.Range('E45').Value = '  -0.93%  '

$ws.Range('E46').Value = '  +2.27%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.007'
$ws.Range('E47').Value = '  +0.23%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '133.89'
$ws.Range('E48').Value = '  +1.56%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.07954'
$ws.Range('E49').Value = '  +0.32%  '

$ws.Range('B50').Value = 'Flow'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.192'
$ws.Range('E50').Value = '  -2.01%  '

$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.186'
$ws.Range('E51').Value = '  -3.18%  '
